$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.076.72'
$ws.Range('E2').Value = '  +0.04%  '
$ws.Range('D3').Value = '1.654.95'
$ws.Range('E3').Value = '  -0.24%  '
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  -0.19%  '
$ws.Range('D5').Value = '214.88'
$ws.Range('E5').Value = '  +3.57%  '
$ws.Range('D6').Value = '0.5249'
$ws.Range('E6').Value = '  +1.60%  '
$ws.Range('D7').Value = '1.001'
$ws.Range('E7').Value = '  -0.17%  '
$ws.Range('D8').Value = '0.2626'
$ws.Range('E8').Value = '  +1.84%  '
$ws.Range('D9').Value = '0.06387'
$ws.Range('D10').Value = '20.82'
$ws.Range('E10').Value = '  -0.20%  '
$ws.Range('D11').Value = '0.07745'
$ws.Range('E11').Value = '  +2.96%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.645.34'
$ws.Range('E12').Value = '  -0.81%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').Value = '4.459'
$ws.Range('E13').Value = '  +1.67%  '
$ws.Range('D14').Value = '1.878.02'
$ws.Range('E14').Value = '  -0.41%  '
$ws.Range('D15').Value = '0.5517'
$ws.Range('E15').Value = '  +2.62%  '
$ws.Range('D16').Value = '0.0₅8311'
$ws.Range('E16').Value = '  +5.32%  '
$ws.Range('D17').Value = '65.13'
$ws.Range('E17').Value = '  -1.31%  '
$ws.Range('D18').Value = '26.098.49'
$ws.Range('E18').Value = '  +0.01%  '
$ws.Range('E19').Value = '  -0.15%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.750'
$ws.Range('E20').Value = '  +1.63%  '
$ws.Range('D21').Value = '190.76'
$ws.Range('E21').Value = '  +2.06%  '
$ws.Range('D22').Value = '10.26'
$ws.Range('E22').Value = '  +1.07%  '
$ws.Range('D23').Value = '6.353'
$ws.Range('E23').Value = '  +2.97%  '
$ws.Range('E24').Value = '  -0.19%  '
$ws.Range('D25').Value = '143.17'
$ws.Range('E25').Value = '  -3.26%  '
$ws.Range('E26').Value = '  +3.69%  '
$ws.Range('D27').Value = '7.415'
$ws.Range('E27').Value = '  +0.80%  '
$ws.Range('D28').Value = '16.03'
$ws.Range('E28').Value = '  +2.80%  '
$ws.Range('D29').Value = '1.419'
$ws.Range('E29').Value = '  +2.81%  '
$ws.Range('D30').Value = '0.05949'
$ws.Range('E30').Value = '  -2.23%  '
$ws.Range('E31').Value = '  +0.19%  '
$ws.Range('D32').Value = '3.443'
$ws.Range('E32').Value = '  -0.32%  '
$ws.Range('D33').Value = '3.415'
$ws.Range('E33').Value = '  +0.84%  '
$ws.Range('D34').Value = '1.658'
$ws.Range('E34').Value = '  +2.16%  '
$ws.Range('D35').Value = '0.9993'
$ws.Range('E35').Value = '  +2.03%  '
$ws.Range('E36').Value = '  +0.62%  '
$ws.Range('D37').Value = '2.764'
$ws.Range('E37').Value = '  +0.60%  '
$ws.Range('D38').Value = '0.5649'
$ws.Range('E38').Value = '  -3.37%  '
$ws.Range('D39').Value = '0.01604'
$ws.Range('E39').Value = '  +0.90%  '
$ws.Range('D40').Value = '5.886'
$ws.Range('E40').Value = '  +0.30%  '
$ws.Range('D41').Value = '0.8567'
$ws.Range('E41').Value = '  +0.92%  '
$ws.Range('E42').Value = '  -0.09%  '
$ws.Range('D43').Value = '1.028.42'
$ws.Range('E43').Value = '  -6.67%  '
$ws.Range('D44').Value = '99.54'
$ws.Range('E44').Value = '  -0.26%  '
$ws.Range('D45').Value = '1.802.46'
$ws.Range('E45').Value = '  -0.56%  '
$ws.Range('E46').Value = '  -2.27%  '
$ws.Range('D47').Value = '55.94'
$ws.Range('E47').Value = '  +2.25%  '
$ws.Range('E48').Value = '  +0.27%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.070'
$ws.Range('E49').Value = '  +0.88%  '
$ws.Range('D50').Value = '0.05155'
$ws.Range('E50').Value = '  -1.36%  '
$ws.Range('D51').Value = '5.988'
$ws.Range('E51').Value = '  +2.48%  '
